$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PythonCode")

# Replace the "search" test-case row (row 5): a new code sample that prints
# via an "output" variable instead of a bare print(), paired with the
# corrected "Submission Successful" label.
$newSearchCode = "def search(input_list, num):`r`nif(num in input_list):`r`noutput = (`"Element Found`")`r`n\b`r`n\b`r`nelse:`r`noutput = (`"Not Found`")`r`n\b`r`n\b`r`nreturn(output)"
$ws.Range("A5").Value2 = $newSearchCode
$ws.Range("B5").Value2 = "Submission Successful"
# Re-fit the row so it keeps its original automatic (non-custom) height
# instead of the explicit height Excel would otherwise stamp on it after
# the multi-line text assignment.
$ws.Range("A5").EntireRow.AutoFit()

# Fix the capitalisation/wording of the "submission success" label used by
# the remaining test-case rows (7, 9, 11) -> "Submission Successful"
$ws.Range("B7").Value2 = "Submission Successful"
$ws.Range("B9").Value2 = "Submission Successful"
$ws.Range("B11").Value2 = "Submission Successful"

# Make the PythonCode sheet the active tab with A5 selected, matching the
# new focus of the workbook.
$ws.Activate()
$ws.Range("A5").Select()
